$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 64.75
$ws.Range("I9").Value = 63.5
$ws.Range("K9").Value = 63.5
$ws.Range("M9").Value = 105.5
$ws.Range("H12").Value = 199.66667
$ws.Range("I12").Value = 199.66667
$ws.Range("K12").Value = 199.66667
$ws.Range("M12").Value = -29.66667000000001
$ws.Range("H18").Value = 1366.6666
$ws.Range("J18").Value = 1899.5
$ws.Range("L18").Value = 1899.5
$ws.Range("N18").Value = -2467.5
$ws.Range("H29").Value = 1003
$ws.Range("I29").Value = 1003
$ws.Range("K29").Value = 3009
$ws.Range("M29").Value = -2728
$ws.Range("H58").Value = 15
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H93").Value = 14480.667
$ws.Range("J93").Value = 14721.5
$ws.Range("L93").Value = 14721.5
$ws.Range("N93").Value = -19713.5
$ws.Range("H138").Value = 3217.3635
$ws.Range("J138").Value = 3327.125
$ws.Range("L138").Value = 9981.375
$ws.Range("N138").Value = -20261.375

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1631.1666
$ws.Range("I88").Value = 1597.6666
$ws.Range("J88").Value = 1664.6666
$ws.Range("K88").Value = 1597.6666
$ws.Range("L88").Value = 1664.6666
$ws.Range("M88").Value = -1191.6666
$ws.Range("N88").Value = -2476.6666
$ws.Range("H91").Value = 1631.1666
$ws.Range("I91").Value = 1597.6666
$ws.Range("J91").Value = 1664.6666
$ws.Range("K91").Value = 1597.6666
$ws.Range("L91").Value = 1664.6666
$ws.Range("M91").Value = -193.6666
$ws.Range("N91").Value = -4472.6666

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 2097.25
$ws.Range("I10").Value = 2199.5
$ws.Range("J10").Value = 1995
$ws.Range("K10").Value = 2199.5
$ws.Range("L10").Value = 1995
$ws.Range("M10").Value = -2059.5
$ws.Range("N10").Value = -2275
$ws.Range("H29").Value = 1500
$ws.Range("I29").Value = 1500
$ws.Range("K29").Value = 1500
$ws.Range("M29").Value = -1211
$ws.Range("H80").Value = 2687.3333
$ws.Range("I80").Value = 2687.5
$ws.Range("K80").Value = 2687.5
$ws.Range("M80").Value = -1689.5
$ws.Range("H83").Value = 2687.3333
$ws.Range("I83").Value = 2687.5
$ws.Range("K83").Value = 13437.5
$ws.Range("M83").Value = -8445.5
$ws.Range("H88").Value = 13163.667
$ws.Range("J88").Value = 13163.667
$ws.Range("L88").Value = 13163.667
$ws.Range("N88").Value = -13975.667
$ws.Range("H91").Value = 13163.667
$ws.Range("J91").Value = 13163.667
$ws.Range("L91").Value = 13163.667
$ws.Range("N91").Value = -15971.667
$ws.Range("H95").Value = 17624
$ws.Range("J95").Value = 17624
$ws.Range("L95").Value = 17624
$ws.Range("N95").Value = -23116

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 1436.6666
$ws.Range("I32").Value = 691
$ws.Range("J32").Value = 2182.3333
$ws.Range("K32").Value = 691
$ws.Range("L32").Value = 2182.3333
$ws.Range("M32").Value = -375
$ws.Range("N32").Value = -2814.3333
$ws.Range("H58").Value = 4357.2856
$ws.Range("I58").Value = 3872
$ws.Range("K58").Value = 3872
$ws.Range("M58").Value = -3669
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").ClearContents()
$ws.Range("N109").Value = 0
$ws.Range("H134").Value = 6099.6665
$ws.Range("I134").Value = 649.5
$ws.Range("K134").Value = 1948.5
$ws.Range("M134").Value = 586.5
$ws.Range("H136").Value = 4357.2856
$ws.Range("I136").Value = 3872
$ws.Range("K136").Value = 11616
$ws.Range("M136").Value = -9066

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1023.6316
$ws.Range("I2").Value = 396.75
$ws.Range("J2").Value = 4367
$ws.Range("K2").Value = 2380.5
$ws.Range("L2").Value = 26202
$ws.Range("M2").Value = -2267.5
$ws.Range("N2").Value = -26428
$ws.Range("H4").Value = 333334
$ws.Range("I4").Value = 500000
$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 1500000
$ws.Range("L4").Value = 6
$ws.Range("M4").Value = -1499888
$ws.Range("N4").Value = -230
$ws.Range("H11").Value = 5000
$ws.Range("J11").Value = 5000
$ws.Range("L11").Value = 15000
$ws.Range("N11").Value = -15280
$ws.Range("H32").Value = 100
$ws.Range("J32").Value = 100
$ws.Range("L32").Value = 300
$ws.Range("N32").Value = -866
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").ClearContents()
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = 0
$ws.Range("H64").Value = 912
$ws.Range("I64").Value = 912
$ws.Range("K64").Value = 2736
$ws.Range("M64").Value = -2466
$ws.Range("H67").Value = 912
$ws.Range("I67").Value = 912
$ws.Range("K67").Value = 2736
$ws.Range("M67").Value = -1800
$ws.Range("H69").Value = 9506
$ws.Range("I69").Value = 9013
$ws.Range("J69").Value = 9999
$ws.Range("K69").Value = 27039
$ws.Range("L69").Value = 29997
$ws.Range("M69").Value = -26228
$ws.Range("N69").Value = -31619
$ws.Range("H72").Value = 9506
$ws.Range("I72").Value = 9013
$ws.Range("J72").Value = 9999
$ws.Range("K72").Value = 81117
$ws.Range("L72").Value = 89991
$ws.Range("M72").Value = -77061
$ws.Range("N72").Value = -98103
$ws.Range("H103").Value = 3399.5
$ws.Range("J103").Value = 4959.4
$ws.Range("L103").Value = 14878.2
$ws.Range("N103").Value = -16636.2
$ws.Range("H109").Value = 2279.6667
$ws.Range("I109").Value = 2279.6667
$ws.Range("K109").Value = 6839.000100000001
$ws.Range("M109").Value = -5799.000100000001
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = 0
$ws.Range("H140").Value = 1426.3334
$ws.Range("I140").Value = 1426.3334
$ws.Range("K140").Value = 4279.0002
$ws.Range("M140").Value = 900.9997999999996

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 3216034.2
$ws.Range("I3").Value = 3335333.2
$ws.Range("J3").Value = 3126559.8
$ws.Range("K3").Value = 3335333.2
$ws.Range("L3").Value = 3126559.8
$ws.Range("M3").Value = -3335217.2
$ws.Range("N3").Value = -3126791.8
$ws.Range("H33").Value = 170
$ws.Range("I33").Value = 170
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 170
$ws.Range("L33").ClearContents()
$ws.Range("N33").Value = 0
$ws.Range("M33").Value = 82
$ws.Range("H70").Value = 7000
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 7000
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6000
$ws.Range("I40").Value = 6000
$ws.Range("K40").Value = 6000
$ws.Range("M40").Value = -5864
$ws.Range("H82").Value = 977.8
$ws.Range("I82").Value = 977.8
$ws.Range("K82").Value = 977.8
$ws.Range("M82").Value = -616.8
$ws.Range("H85").Value = 977.8
$ws.Range("I85").Value = 977.8
$ws.Range("K85").Value = 977.8
$ws.Range("M85").Value = 270.2
$ws.Range("H108").Value = 30000
$ws.Range("J108").Value = 30000
$ws.Range("L108").Value = 30000
$ws.Range("N108").Value = -37680
$ws.Range("H136").Value = 1443.6666
$ws.Range("I136").Value = 892.4
$ws.Range("K136").Value = 2677.2
$ws.Range("M136").Value = -127.1999999999998

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4998.6665
$ws.Range("H65").Value = 4998.6665
$ws.Range("H132").Value = 4002
$ws.Range("I132").Value = 4002
$ws.Range("K132").Value = 12006
$ws.Range("M132").Value = -9476
